$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 270, shifting existing row 270 (IAD) and below down by one.
$ws.Rows.Item(270).Insert()

# Copy style from row 271 (the old row 270, now shifted down, which has the colo-code bold style in col A)
$ws.Cells.Item(271, 1).Copy()
$ws.Cells.Item(270, 1).PasteSpecial(-4122)

# Populate the new row with the Malang, Indonesia colo entry.
$ws.Cells.Item(270, 1).Value = "MLG"
$ws.Cells.Item(270, 2).Value = "Malang, Indonesia"
$ws.Cells.Item(270, 3).Value = "Asia Pacific"
$ws.Cells.Item(270, 4).Value = "Malang"
$ws.Cells.Item(270, 5).Value = "Indonesia"
$ws.Cells.Item(270, 6).Value = "ID"
$ws.Cells.Item(270, 7).Value = -8.100346999999999
$ws.Cells.Item(270, 8).Value = 112.186641
